# ----------------------------------------------------------------------
# Gantt_Diagram_Robot_Project_2018.xlsx  -  "Add Gant Gyroscope stuff"
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")

# ------------------------------------------------------------------
# 1) Extend the Gantt grid by 4 weeks (weeks 9-12, columns BO:CP)
#    by first cloning the formatting of the previous 4-week block
#    (columns AM:BN) into the new block.
# ------------------------------------------------------------------
$srcBlock = $ws.Range("AM1:BN42")
$srcBlock.Copy()
$ws.Range("BO1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Column widths for the newly used columns (match neighbouring week style)
$ws.Range("BO1:CP1").EntireColumn.ColumnWidth = 2.42578125

# ------------------------------------------------------------------
# 2) Row 6 : running date chain (each cell = previous cell + 1)
# ------------------------------------------------------------------
$row6cols = @("BO","BP","BQ","BR","BS","BT","BU","BV","BW","BX","BY","BZ","CA","CB","CC","CD","CE","CF","CG","CH","CI","CJ","CK","CL","CM","CN","CO","CP")
$prevCol = "BN"
foreach ($col in $row6cols) {
    $ws.Range($col + "6").Formula = "=" + $prevCol + "6+1"
    $prevCol = $col
}

# ------------------------------------------------------------------
# 3) Row 7 : weekday letter for each day, based on row 6 date
# ------------------------------------------------------------------
foreach ($col in $row6cols) {
    $ws.Range($col + "7").Formula = "=CHOOSE(WEEKDAY(" + $col + "6,1),""S"",""M"",""T"",""W"",""T"",""F"",""S"")"
}

$excel.Calculate()

# ------------------------------------------------------------------
# 4) Row 4 : "Week N" header (merged, one per week) + row 5 start date
# ------------------------------------------------------------------
$weekBlocks = @(
    @{Start="BO"; End="BU"},
    @{Start="BV"; End="CB"},
    @{Start="CC"; End="CI"},
    @{Start="CJ"; End="CP"}
)

foreach ($wk in $weekBlocks) {
    $startCol = $wk.Start
    $endCol   = $wk.End

    $ws.Range($startCol + "4").Formula = "=""Week ""&(" + $startCol + "6-(`$C`$4-WEEKDAY(`$C`$4,1)+2))/7+1"
    $ws.Range($startCol + "4:" + $endCol + "4").Merge()

    $ws.Range($startCol + "5").Formula = "=" + $startCol + "6"
    $ws.Range($startCol + "5:" + $endCol + "5").Merge()
}

$excel.Calculate()

# ------------------------------------------------------------------
# 5) Conditional formatting ("today" highlight) for the new weeks,
#    mirroring the template used for the existing weeks (K6:BN7 /
#    K6:BN41): a red border rule and a filled/white-font rule.
# ------------------------------------------------------------------
foreach ($wk in $weekBlocks) {
    $startCol = $wk.Start
    $endCol   = $wk.End
    $rng = $ws.Range($startCol + "6:" + $endCol + "7")

    $fc1 = $rng.FormatConditions.Add(2, 0, "=" + $startCol + "`$6=TODAY()")
    $fc1.Borders.Item(7).LineStyle = 1
    $fc1.Borders.Item(7).Color = 192
    $fc1.Borders.Item(10).LineStyle = 1
    $fc1.Borders.Item(10).Color = 192

    $fc2 = $rng.FormatConditions.Add(2, 0, "=" + $startCol + "`$6=TODAY()")
    $fc2.Font.ThemeColor = 1
    $fc2.Interior.ThemeColor = 6
}

# ------------------------------------------------------------------
# 6) Task rows 22-27 : new "Gyroscope/Accelerometer" task group
# ------------------------------------------------------------------
$ws.Range("B22").Value = "Gyroscope/Accelerometer"

$ws.Range("B23").Value = "I2C protocol"
$ws.Range("C23").Value = "M,Guillaume Nicolle"
$ws.Range("E23").Value = 43384
$ws.Range("G23").Value = 25

$ws.Range("B24").Value = "Data Treatment"
$ws.Range("C24").Value = "M,Guillaume Nicolle"
$ws.Range("E24").Value = 43384
$ws.Range("G24").Value = 2

$ws.Range("B25").Value = "Integration to Amine's program"
$ws.Range("C25").Value = "M,Guillaume Nicolle"
$ws.Range("E25").Value = 43384
$ws.Range("G25").Value = 7

$ws.Range("B26").Value = "Noise treatment"
$ws.Range("C26").Value = "M,Guillaume Nicolle"
$ws.Range("E26").Value = 43384
$ws.Range("G26").Value = 14

$ws.Range("B27").Value = "Wiki Page"
$ws.Range("C27").Value = "M,Guillaume Nicolle"
$ws.Range("E27").Value = 43384
$ws.Range("G27").Value = 1

$excel.Calculate()

# ------------------------------------------------------------------
# 7) Row 11 height tweak
# ------------------------------------------------------------------
$ws.Range("B11").EntireRow.RowHeight = 18

# ------------------------------------------------------------------
# 8) View state: frozen-pane top-left cell and active selection
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$ws.Range("BO23").Select()

$excel.Calculate()
